$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E5").Value = $false
$ws.Range("E6").Value = $false

$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "ice cream"
$ws.Range("C8").Value = 10
$ws.Range("D8").Value = $false
$ws.Range("E8").Value = $true
